# Insert a new record at row 21 (shifts the existing rows 21-48 down to 22-49)
# and populate it with the new Cilantro price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44792
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 100112040
$ws.Range("G21").Value = "Cilantro"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 12000
$ws.Range("N21").Value = "`$/caja 36 atados"
$ws.Range("O21").Value = "Provincia de Quillota"
$ws.Range("P21").Value = 333
$ws.Range("Q21").Value = 36
$ws.Range("R21").Value = "Hortaliza"
